# MarcheTemplate.xlsx edit:
# - Swap the two placeholder/date-format labels shown under the header row
#   (E2 becomes the "YYYY-MM-DD" label, B2 becomes the "****/**" placeholder).
#   Leading "'" keeps them as quote-prefixed text cells, same as the source file.
# - Move the active selection to B3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set E2 first so the new shared-string entries are appended in the same
# order as the target workbook (YYYY-MM-DD before ****/**).
$ws.Range("E2").Value = "'YYYY-MM-DD"
$ws.Range("B2").Value = "'****/**"

$ws.Range("B3").Select()
